$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item("TextBox 3")
$tr = $shape.TextFrame.TextRange

# The caption textbox currently stores its text as a series of separate
# runs ("Followed" / " " / "by" / " " / "a" / " " / "picture"). Consolidate
# them into a single run, like a PowerPoint writer that slims down output
# by avoiding unnecessary <a:r> elements would.
$tr.Delete()
$tr.Text = "Followed by a picture"
